$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 770, shifting existing rows 770-811 down to 771-812
$ws.Rows("770:770").Insert()

# Column A holds dates stored as plain text (e.g. "2026/02/07"), not real Excel
# dates. Force text format before assignment so Excel doesn't auto-convert the
# string into a date serial number, then clear the explicit formatting so the
# cell keeps default (unstyled) formatting like the rest of the data rows.
$ws.Range("A770").NumberFormat = "@"
$ws.Range("A770").Value() = "2026/02/07"
$ws.Range("A770").ClearFormats()

$ws.Range("B770").Value() = "土"
$ws.Range("C770").Value() = 19
$ws.Range("D770").Value() = 90
